$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string reorderings (country name ranking swaps) ---
# Cuba, Ruanda, Surinam, Jamaica -> Cuba, Surinam, Ruanda, Jamaica
$ws.Range("A123").Value = "Surinam"
$ws.Range("A124").Value = "Ruanda"

# Tailandia, Gambia, Somalia, Angola, Lituania -> Tailandia, Angola, Gambia, Somalia, Lituania
$ws.Range("A130").Value = "Angola"
$ws.Range("A131").Value = "Gambia"
$ws.Range("A132").Value = "Somalia"

# Niger, Chad, Vietnam, Polinesia Francesa, Martinica -> Niger, Polinesia Francesa, Chad, Vietnam, Martinica
$ws.Range("A165").Value = "Polinesia Francesa"
$ws.Range("A166").Value = "Republica del Chad"
$ws.Range("A167").Value = "Vietnam"

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 02:39"

# --- Numeric data updates ---
$ws.Range("B4").Value = 6748842
$ws.Range("C4").Value = 37772
$ws.Range("D4").Value = 4024385
$ws.Range("E4").Value = 2525483
$ws.Range("G4").Value = 454
$ws.Range("H4").Value = 198974
$ws.Range("B13").Value = 565446
$ws.Range("C13").Value = 9909
$ws.Range("E13").Value = 124826
$ws.Range("G13").Value = 315
$ws.Range("H13").Value = 11667
$ws.Range("B108").Value = 7244
$ws.Range("C108").Value = 6
$ws.Range("D108").Value = 6555
$ws.Range("E108").Value = 565
$ws.Range("B114").Value = 5104
$ws.Range("C114").Value = 29
$ws.Range("D114").Value = 4374
$ws.Range("E114").Value = 629
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 101
$ws.Range("B118").Value = 4934
$ws.Range("C118").Value = 6
$ws.Range("E118").Value = 959
$ws.Range("B121").Value = 4772
$ws.Range("C121").Value = 23
$ws.Range("D121").Value = 1828
$ws.Range("E121").Value = 2882
$ws.Range("B123").Value = 4611
$ws.Range("C123").Value = 29
$ws.Range("D123").Value = 3935
$ws.Range("E123").Value = 581
$ws.Range("G123").Value = 2
$ws.Range("H123").Value = 95
$ws.Range("B124").Value = 4602
$ws.Range("C124").Value = 11
$ws.Range("D124").Value = 2736
$ws.Range("E124").Value = 1844
$ws.Range("H124").Value = 22
$ws.Range("B125").Value = 3933
$ws.Range("C125").Value = 162
$ws.Range("D125").Value = 1161
$ws.Range("E125").Value = 2728
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 44
$ws.Range("B130").Value = 3439
$ws.Range("C130").Value = 51
$ws.Range("D130").Value = 1324
$ws.Range("E130").Value = 1979
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = 136
$ws.Range("B131").Value = 3405
$ws.Range("D131").Value = 1723
$ws.Range("E131").Value = 1579
$ws.Range("H131").Value = 103
$ws.Range("B132").Value = 3389
$ws.Range("D132").Value = 2803
$ws.Range("E132").Value = 488
$ws.Range("H132").Value = 98
$ws.Range("B158").Value = 1534
$ws.Range("C158").Value = 8
$ws.Range("D158").Value = 1282
$ws.Range("E158").Value = 230
$ws.Range("B165").Value = 1099
$ws.Range("C165").Value = 146
$ws.Range("D165").Value = 672
$ws.Range("E165").Value = 425
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 2
$ws.Range("B166").Value = 1085
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 940
$ws.Range("E166").Value = 64
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 81
$ws.Range("B167").Value = 1063
$ws.Range("D167").Value = 926
$ws.Range("E167").Value = 102
$ws.Range("H167").Value = 35
